$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# New note texts (appended catatan revisi)
$notes = @(
    '[4/28/2015, 03:58] Pak Devi Leuwigajah: Za yang artesis menu nya belum ada paling di tambahin jumlah pengguna aja klo artesis',
    '[4/28/2015, 03:58] Pak Devi Leuwigajah: Blm ada isi nya masih error 404',
    '[4/28/2015, 04:01] Pak Devi Leuwigajah: Terus untuk jalan ga usah ada ketersediaan lahan za',
    '[4/28/2015, 04:06] Pak Devi Leuwigajah: Dan untuk kategori jalan lebar 1-2m= jalan setapak. 2-4=jalan lingkungan. 4-12= jalan utama',
    '[4/28/2015, 04:07] Pak Devi Leuwigajah: Atau untuk ketersediaan lahan untuk jalan isi nya cuma ada atau tidak aja jangan pake angka',
    '[4/28/2015, 04:42] Pak Devi Leuwigajah: Za terus untuk posisi koordinat yang longtitude latitude nya di seragamin aja',
    '[4/28/2015, 04:43] Pak Devi Leuwigajah: Jadi longtitude awal latitude awal terus longtitude akhir latitude akhir soal na jalan sama drainase beda beda bisi ke lieur entri data na',
    '[4/28/2015, 04:59] Pak Devi Leuwigajah: Terus untuk yang "sedang dilaksanakan" kata kata nya di tambah za "sedang / akan dilaksanakan "'
)

$startRow = 79
$lastRow = $startRow + $notes.Count - 1   # 86

# Copy formatting of the last existing data row (plain style, no border) down to
# the rows that will hold the 7 new non-final notes.
$ws.Range("A77:B77").Copy($ws.Range("A$startRow`:B$($lastRow-1)"))

# Copy formatting of the current final row (bottom-border style) to the new final row.
$ws.Range("A78:B78").Copy($ws.Range("A$lastRow`:B$lastRow"))

# Fill in the text values for the new rows.
for ($i = 0; $i -lt $notes.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value2 = $notes[$i]
}

# Widen column A slightly to accommodate the longer notes (matches the resulting
# manual column width after the edit; 154.8 is the closest input that this
# runtime's pixel-quantized width model resolves to the target 155.71 chars).
$ws.Columns.Item(1).ColumnWidth = 154.8

# Grow the table / autofilter to cover the newly added rows.
$lo.Resize($ws.Range("A1:B$lastRow"))

# Restore the active selection near the newly added content.
$ws.Range("A81").Select()
